# Automatische test-sync: 2025-08-04 20:49:50
# Appends a new incoming-mail log row (row 23) to the "Logs" sheet,
# extends the conditional-formatting ranges to cover it, and bumps the
# "Retour / Terugbetaling" tally on the "Dashboard" sheet from 2 to 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 23

$ws.Cells.Item($newRow, 1).Value = "Mijn retour is nog steeds niet verwerkt."
$ws.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 3).Value = "Testmail #11: Mijn retour is nog steeds niet verwerkt."
$ws.Cells.Item($newRow, 4).Value = "Retour / Terugbetaling"
$ws.Cells.Item($newRow, 5).Value = "Beste klant,`nDank u voor uw bericht. Om uw retourzending verder te kunnen onderzoeken, heb ik wat meer informatie nodig. Kunt u alstublieft uw ordernummer en de datum van de retourzending doorgeven? Op die manier kunnen wij u sneller van dienst zijn.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$ws.Cells.Item($newRow, 6).Value = "2025-08-04 20:49:17"
$ws.Cells.Item($newRow, 7).Value = "Ja"
$ws.Cells.Item($newRow, 8).Value = "Nee"
$ws.Cells.Item($newRow, 9).Value = "Ja"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# The multi-line "Antwoord" text would otherwise leave the new row with a
# stale auto-estimated custom height; re-fit it back down to the sheet's
# standard row height so row 23 matches its siblings (no customHeight).
$ws.Rows.Item($newRow).AutoFit()

# Extend the conditional-formatting coverage for each column from row 22
# down to the newly-added row 23 (one ModifyAppliesToRange per group is
# enough -- every cfRule sharing that <conditionalFormatting> block moves
# together).
$ws.Range("D2:D22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D23"))
$ws.Range("G2:G22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G23"))
$ws.Range("H2:H22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H23"))
$ws.Range("I2:I22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I23"))
$ws.Range("J2:J22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J23"))

# Update the Dashboard summary count for "Retour / Terugbetaling" (2 -> 3).
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B5").Value = 3
